$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new daily record (2020-08-12) as row 74, right after the
# existing last data row (73).
$newRow = 74

# Column A holds a date-like label ("2020-08-12") but the sheet stores it
# as plain text (shared string), not an Excel date serial. Force text
# entry by switching the cell to a text number format before assigning
# the value, then restore the default/general format so the cell keeps
# no explicit style (matching the rest of the data rows).
$dateCell = $ws.Cells.Item($newRow, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2020-08-12"
$dateCell.ClearFormats()

$ws.Cells.Item($newRow, 2).Value = 498380
$ws.Cells.Item($newRow, 3).Value = 545262
$ws.Cells.Item($newRow, 4).Value = 83473
$ws.Cells.Item($newRow, 5).Value = 54666
$ws.Cells.Item($newRow, 6).Value = 26.49
